# Remove four bullet items from the bug list:
#   - "Player laser spawn place moves up every game"
#   - "Player is an explosion"
#   - "Enemies do not damage the shields"
#   - "Hits side wall and ends game" (including its bookmark)
#
# Delete from the bottom up so earlier paragraph indices remain valid
# as later paragraphs are removed.

$d = $word.ActiveDocument

$d.Paragraphs.Item(8).Range.Delete()  # "Hits side wall and ends game"
$d.Paragraphs.Item(4).Range.Delete()  # "Enemies do not damage the shields"
$d.Paragraphs.Item(3).Range.Delete()  # "Player is an explosion"
$d.Paragraphs.Item(2).Range.Delete()  # "Player laser spawn place moves up every game"
